# Insert a new price record as the new row 223, pushing the existing
# rows 223-265 down to 224-266 (weekly Fruit/Hortaliza price update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(223).Insert()

$ws.Range("A223").Value = 4
$ws.Range("B223").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C223").Value = "Los Lagos"
$ws.Range("D223").Value = 44711
$ws.Range("E223").Value = 10
$ws.Range("F223").Value = 100112037
$ws.Range("G223").Value = "Cebollín"
$ws.Range("H223").Value = "Sin especificar"
$ws.Range("I223").Value = "Primera"
$ws.Range("J223").Value = 60
$ws.Range("K223").Value = 11000
$ws.Range("L223").Value = 12000
$ws.Range("M223").Value = 11500
$ws.Range("N223").Value = "$/paquete 36 unidades"
$ws.Range("O223").Value = "Región Metropolitana"
$ws.Range("P223").Value = 319
$ws.Range("Q223").Value = 36
$ws.Range("R223").Value = "Hortaliza"
